# Update automatico via Actualizar 06-17-2020 05-26-11
# Adds the 2020-06-16 (serial 43998) row of new COVID-19 patient-condition
# data to the "Condicion_Pacientes" table on Hoja1, growing the table/
# worksheet from A1:F95 to A1:F96.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Grow the table by one row; Excel will extend the table range (and the
# AutoFilter range) from A1:F95 to A1:F96 automatically.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$newRow = $lo.ListRows.Add()

# Copy the formatting of the previous last data row (row 95) down into the
# freshly added row 96 so the new cells pick up the same number formats /
# alignment styles used throughout the table.
$ws.Range("A95:F95").Copy() | Out-Null
$ws.Range("A96:F96").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the new row's values.
$ws.Cells.Item(96, 1).Value = 43998
$ws.Cells.Item(96, 2).Value = 1086
$ws.Cells.Item(96, 3).Value = 478
$ws.Cells.Item(96, 4).Value = 568
$ws.Cells.Item(96, 5).Value = 302
$ws.Cells.Item(96, 6).Value = 51

# Match the author's resulting selection (cell A96 active).
$ws.Range("A96").Select() | Out-Null
